# fix: issue in author names
# Remove duplicated author-name entries that had been accidentally
# repeated within the "Authors" column (E2, E3, E4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Saleh A Alqahtani, `nJörn M Schattenberg"

$ws.Range("E3").Value = "Zhi-Yong Li, `nTuya, `nHai-Tao Li, `nJiang He, `nQuesheng, `nGuang-Ping Dong, `nMing-Shuo Zhang, `nJian-Qin Liu, `nXiu-Lan Huang, `nXiao-Rong Wang, `nMakabel Bolat, `nXin Feng, `nFang-Bo Zhang, `nFeng Jiang"

$ws.Range("E4").Value = "James R Vallerand, `nRyan E Rhodes, `nGordan J Walker, `nKerry S Courneya"
